# Applies the commit "[base commands] - [assertMatch(text,regex)]: NEW command to
# check for text value via regular expression" to the hidden '#system' lookup
# sheet that backs the dependent dropdowns on the ImportCSV sheet.
#
# Net effect (mirrors the authoritative xlsx diff):
#   1. Insert "assertMatch(text,regex)" into the alphabetically sorted 'base'
#      command list (column F), pushing the following entries down a row.
#   2. Insert "openFile(filePath)" into the alphabetically sorted 'external'
#      command list (column J), pushing the following entries down a row.
#   3. Remove the 'tn.5250' category altogether: its entry disappears from the
#      'target' category list (column A) and its whole data column (AA) is
#      deleted, so columns AB:AG (web, webalert, webcookie, ws, ws.async, xml)
#      shift one column to the left (AA:AF).
#   4. Keep the named ranges in sync with the new extents.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

function Col-Letter([int]$n) {
    $s = ""
    while ($n -gt 0) {
        $rem = ($n - 1) % 26
        $s = [char](65 + $rem) + $s
        $n = [int](($n - 1) / 26)
    }
    return $s
}

# ---------------------------------------------------------------------------
# 1) "base" (column F = 6): insert new row at F11, shifting F11:F44 -> F12:F45
# ---------------------------------------------------------------------------
$baseCol = 6
for ($r = 44; $r -ge 11; $r--) {
    $src = $ws.Cells.Item($r, $baseCol).Value2
    $ws.Cells.Item($r + 1, $baseCol).Value = $src
}
$ws.Cells.Item(11, $baseCol).Value = "assertMatch(text,regex)"

# ---------------------------------------------------------------------------
# 2) "external" (column J = 10): insert new row at J2, shifting J2:J6 -> J3:J7
# ---------------------------------------------------------------------------
$externalCol = 10
for ($r = 6; $r -ge 2; $r--) {
    $src = $ws.Cells.Item($r, $externalCol).Value2
    $ws.Cells.Item($r + 1, $externalCol).Value = $src
}
$ws.Cells.Item(2, $externalCol).Value = "openFile(filePath)"

# ---------------------------------------------------------------------------
# 3) Delete the whole "tn.5250" data column (AA = 27), shifting AB:AG (28:33)
#    left into AA:AF (27:32) for every row used by the sheet (1-151).
# ---------------------------------------------------------------------------
$lastDataRow = 151
$deletedCol = 27
$lastUsedCol = 33
for ($r = 1; $r -le $lastDataRow; $r++) {
    for ($c = $deletedCol; $c -lt $lastUsedCol; $c++) {
        $src = $ws.Cells.Item($r, $c + 1).Value2
        $ws.Cells.Item($r, $c).Value = $src
    }
    $ws.Cells.Item($r, $lastUsedCol).ClearContents()
}

# ---------------------------------------------------------------------------
# 4) Remove "tn.5250" from the "target" category list (column A = 1): delete
#    A27, shifting A28:A33 up into A27:A32.
# ---------------------------------------------------------------------------
$targetCol = 1
$tnRow = 27
$lastTargetRow = 33
for ($r = $tnRow; $r -lt $lastTargetRow; $r++) {
    $src = $ws.Cells.Item($r + 1, $targetCol).Value2
    $ws.Cells.Item($r, $targetCol).Value = $src
}
$ws.Cells.Item($lastTargetRow, $targetCol).ClearContents()

# ---------------------------------------------------------------------------
# 5) Keep defined names in sync with the new ranges.
# ---------------------------------------------------------------------------
$wb.Names.Item("base").RefersTo = "='#system'!`$F`$2:`$F`$45"
$wb.Names.Item("external").RefersTo = "='#system'!`$J`$2:`$J`$7"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$32"
$wb.Names.Item("web").RefersTo = "='#system'!`$AA`$2:`$AA`$151"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$AB`$2:`$AB`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AC`$2:`$AC`$10"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AD`$2:`$AD`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AE`$2:`$AE`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AF`$2:`$AF`$27"

Write-Output "done"
